$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlUp = -4162
$lastRow = $ws.Cells.Item(1048576, 3).End($xlUp).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    $cell.Value = $current + 1
}
